$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column C (Duration), shifting it to D.
$ws.Range("C1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "Chord"

# Fill the new "Chord" column with boolean FALSE for every data row (rows 2-84).
$lastRow = 84
$chordRange = $ws.Range("C2:C$lastRow")
$chordRange.Value = $false
